# Automatische test-sync: 2025-06-26 19:30:50
# Adds a new log row (row 10) to the "Logs" sheet, extends the
# conditional-formatting ranges to include it, and bumps the
# "Bestelling / Levering" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet -----------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A10").Value = "Bestel je 5 paar werkhandschoenen voor me"
$logs.Range("B10").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Range("C10").Value = "He Johan`nzou je 5 paar werkhandschoenen voor mij kunnen bestellen?`nMarc`nSent using {0}"
$logs.Range("D10").Value = "Bestelling / Levering"
$logs.Range("E10").Value = "Beste Marc,`nBedankt voor je e-mail. Natuurlijk kan ik 5 paar werkhandschoenen voor je bestellen. Kun je me wat meer details geven, zoals het type werkhandschoenen en de maat die je nodig hebt? Als je specifieke voorkeuren hebt, laat het me dan weten zodat ik de bestelling correct kan plaatsen.`nAlvast bedankt!`nMet vriendelijke groet,`nJohan"
$logs.Range("F10").Value = "2025-06-26 19:30:33"
$logs.Range("G10").Value = "Ja"
$logs.Range("H10").Value = "Nee"
$logs.Range("I10").Value = "Ja"

# Cells C10/E10 contain embedded line breaks, which makes the engine pin an
# explicit row height (ht/customHeight) on write. Auto-fitting the row clears
# that override again, matching the rest of the sheet (no explicit row
# heights anywhere else).
$logs.Rows.Item(10).AutoFit()

# Extend the existing conditional formatting ranges (D/G/H/I 2:9 -> 2:10)
# by re-pointing each existing rule at the enlarged range instead of
# creating new duplicate rules.
$colsToExtend = @("D", "G", "H", "I")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range("$($col)2:$($col)9")
    $newRange = $logs.Range("$($col)2:$($col)10")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet --------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 5
